$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B27: change from text "3" to a real numeric value 3
$ws.Range("B27").Value = 3

# Add new row 28 with the new annotation data
$ws.Range("A28").Value = "Ruilin"

# B28 must stay text "4" (not get auto-converted to a number)
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "4"
$ws.Range("B28").Style = "Normal"

$ws.Range("C28").Value = "is appealing and an important open problem."
$ws.Range("D28").Value = "APC"
$ws.Range("E28").Value = "OTH"
$ws.Range("F28").Value = "f913699b-da49-47c6-8043-88c593733ae2"
$ws.Range("G28").Value = "BJyy3a0Ez_annotated.xlsx"
$ws.Range("H28").Value = "The idea of model-parallelism (as opposed to data parallelism) is appealing and an important open problem."
